$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Publish")

# Insert a new column before column G ("ModernGrpLoc" and its data shift right to H)
$ws.Range("G1").EntireColumn.Insert()

# New column header + bold style to match the other header cells
$ws.Range("G1").Value = "OneDriveLoc"
$ws.Range("G1").Font.Bold = $true

# New column data
$ws.Range("G2").Value = "All"
$ws.Range("G3").Value = "All"

# Autofit the new column to its content (matches bestFit width behavior)
$ws.Columns.Item(7).AutoFit()

# The hyperlink anchored at the old G3 now lives at H3 - rebuild it there.
# The source keeps the cell's own text ("ALYAMG-...") distinct from the
# hyperlink's display/tooltip text ("ALYAOG-..."), so restore the cell text
# after Excel's Hyperlinks.Add overwrites it with the display text.
$existingText = $ws.Range("H3").Value2
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("H3"), "mailto:ALYAOG-ADM-AlleExternen@alyaconsulting031.onmicrosoft.com", [Type]::Missing, [Type]::Missing, "ALYAOG-ADM-AlleExternen@alyaconsulting031.onmicrosoft.com")
$ws.Range("H3").Value = $existingText

# Source workbook keeps the hyperlink cell in the default (non-hyperlink) style,
# so undo the auto-applied "Hyperlink" style just like the original had.
$ws.Range("H3").ClearFormats()
